$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph using Find.
$findRange = $d.Content
$found = $findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Determine which paragraph (1-based index) contains the found text.
    $idx = 0
    $targetIdx = -1
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if (($p.Range.Start -le $findRange.Start) -and ($p.Range.End -ge $findRange.End)) {
            $targetIdx = $idx
        }
    }

    # Remove the blank paragraph immediately before it, the "Ver no Jupiter..."
    # paragraph itself, and the "© 2020 ..." paragraph immediately after it, by
    # deleting the single range spanning from the start of the blank paragraph
    # through the end of the copyright paragraph (including its paragraph mark).
    $prevPara = $d.Paragraphs.Item($targetIdx - 1)
    $nextPara = $d.Paragraphs.Item($targetIdx + 1)

    $delRange = $d.Range($prevPara.Range.Start, $nextPara.Range.End)
    $delRange.Delete()
}
